# Add a new quarterly column (BB) to the forecast sheet, mirroring the
# existing pattern: header date in row 1, and forecast values in rows 3-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date for column BB (copy style/format from BA1, then set value)
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# New forecast values for column BB, rows 3-18 repeat the last known value
# from column BA (carried forward), rows 19-21 contain updated forecasts.
$ws.Range("BB3").Value = 3.311198277644523
$ws.Range("BB4").Value = 2.295764853811622
$ws.Range("BB5").Value = 0.9560359274609631
$ws.Range("BB6").Value = 1.33431011236409
$ws.Range("BB7").Value = 0.8625324730765715
$ws.Range("BB8").Value = 1.127663699749437
$ws.Range("BB9").Value = 2.036722793063217
$ws.Range("BB10").Value = 4.211110015923047
$ws.Range("BB11").Value = 1.551589743072856
$ws.Range("BB12").Value = 1.455651414914461
$ws.Range("BB13").Value = 2.211901852239651
$ws.Range("BB14").Value = 3.454003806423267
$ws.Range("BB15").Value = 2.895789704837681
$ws.Range("BB16").Value = 0.486919784128137
$ws.Range("BB17").Value = -0.8339987268308979
$ws.Range("BB18").Value = 2.721554796335779
$ws.Range("BB19").Value = 2.622852459381209
$ws.Range("BB20").Value = 2.447176337618551
$ws.Range("BB21").Value = 2.720226860657204
